# Implement toEven and toOdd operators in the sampleMath workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New sample rows: EVEN / EVEN(negative) / ODD / ODD(negative) ---
$ws.Range("E13").Value = "EVEN"
$ws.Range("F13").Formula = "=EVEN(2.91)"

$ws.Range("E14").Value = "EVEN -"
$ws.Range("F14").Formula = "=EVEN(-7.8)"

$ws.Range("E15").Value = "ODD"
$ws.Range("F15").Formula = "=ODD(2.91)"

$ws.Range("E16").Value = "ODD -"
$ws.Range("F16").Formula = "=ODD(-7.8)"

# --- Defined names for the new sample cells ---
$wb.Names.Add("EvenVal", $ws.Range("F13"))
$wb.Names.Add("EvenNegVal", $ws.Range("F14"))
$wb.Names.Add("OddVal", $ws.Range("F15"))
$wb.Names.Add("OddNegVal", $ws.Range("F16"))

# --- Update the view: move the selection (also clears the stale topLeftCell scroll) ---
$ws.Range("F15").Select() | Out-Null
